$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 247 ("Fruta / hortaliza, semanal").
# This pushes the former rows 247-303 down to 248-304 and the sheet's used
# range grows from A1:R303 to A1:R304.
$ws.Rows(247).Insert()

$ws.Range("A247").Value = 6
$ws.Range("B247").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C247").Value = "Metropolitana"
$ws.Range("D247").Value = 44511
$ws.Range("E247").Value = 13
$ws.Range("F247").Value = 100112043
$ws.Range("G247").Value = "Pepino ensalada"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 1300
$ws.Range("K247").Value = 5000
$ws.Range("L247").Value = 6000
$ws.Range("M247").Value = 5577
$ws.Range("N247").Value = "`$/caja 50 unidades"
$ws.Range("O247").Value = "Región de Arica y Parinacota"
$ws.Range("P247").Value = 112
$ws.Range("Q247").Value = 50
$ws.Range("R247").Value = "Hortaliza"
